# "Generate Report for Archive"
#
# The localization status for the two handed-off source files moved from
# "Ready for handoff" to "In Translation". That status string appears in
# three worksheets:
#   - Overview : columns E (zh-cn) and F (de-de), rows 2-3
#   - zh-cn    : column C (Status), rows 2-3
#   - de-de    : column C (Status), rows 2-3
#
# Because the status text got shorter, the report generator that produced
# this workbook also re-sized (auto-fit) the columns holding that text, so
# we shrink the same columns to match.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Overview sheet --------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# Columns E & F held the status text that just got shorter - narrow them
# to the new auto-fit width.
$wsOverview.Columns.Item(5).ColumnWidth = 12.576851254417766
$wsOverview.Columns.Item(6).ColumnWidth = 12.576851254417766

# --- zh-cn sheet -------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsZhCn.Columns.Item(3).ColumnWidth = 12.576851254417766

# --- de-de sheet -------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

$wsDeDe.Columns.Item(3).ColumnWidth = 12.576851254417766
